$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '86.631.54'
$ws.Range("E2").Value = '  -3.27%  '
$ws.Range("D3").Value = '3.136.23'
$ws.Range("E3").Value = '  -6.98%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '203.95'
$ws.Range("E5").Value = '  -7.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '604.05'
$ws.Range("E6").Value = '  -7.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.364'
$ws.Range("E7").Value = '  -9.75%  '
$ws.Range("E8").Value = '  +7.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.998'
$ws.Range("D10").Value = '3.129.52'
$ws.Range("E10").Value = '  -6.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.525'
$ws.Range("E11").Value = '  -11.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.175'
$ws.Range("E12").Value = '  +4.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000238'
$ws.Range("E13").Value = '  -17.89%  '
$ws.Range("B14").Value = 'Toncoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.21'
$ws.Range("E14").Value = '  -5.94%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.697.42'
$ws.Range("E15").Value = '  -6.61%  '
$ws.Range("D16").Value = '86.270.45'
$ws.Range("E16").Value = '  -3.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '31.61'
$ws.Range("E17").Value = '  -14.10%  '
$ws.Range("D18").Value = '3.124.40'
$ws.Range("E18").Value = '  -6.46%  '
$ws.Range("B19").Value = 'SuiNetwork'
$ws.Range("C19").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.92'
$ws.Range("E19").Value = '  -7.42%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.21'
$ws.Range("E20").Value = '  -11.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '407.66'
$ws.Range("E21").Value = '  -11.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.36'
$ws.Range("E22").Value = '  -12.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.01'
$ws.Range("E23").Value = '  -9.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.07'
$ws.Range("E24").Value = '  -8.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.63'
$ws.Range("E25").Value = '  -10.18%  '
$ws.Range("D26").Value = '3.299.06'
$ws.Range("E26").Value = '  -5.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '72.49'
$ws.Range("E27").Value = '  -7.98%  '
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000127'
$ws.Range("E28").Value = '  -10.69%  '
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  -23.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.996'
$ws.Range("E31").Value = '  -0.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '527.99'
$ws.Range("E32").Value = '  -11.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.16'
$ws.Range("E33").Value = '  -12.49%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.84'
$ws.Range("E34").Value = '  -13.74%  '
$ws.Range("E35").Value = '  -20.89%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.45'
$ws.Range("E36").Value = '  -12.76%  '
$ws.Range("E37").Value = '  -9.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '21.42'
$ws.Range("E38").Value = '  -8.57%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.998'
$ws.Range("E39").Value = '  +0.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '21.73'
$ws.Range("E40").Value = '  -0.64%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.94'
$ws.Range("E41").Value = '  -7.90%  '
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.87'
$ws.Range("E43").Value = '  -14.56%  '
$ws.Range("B44").Value = 'PolygonEcosystemToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.365'
$ws.Range("E44").Value = '  -14.82%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '148.43'
$ws.Range("E45").Value = '  -5.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '170.06'
$ws.Range("E46").Value = '  -10.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '42.71'
$ws.Range("E47").Value = '  -7.68%  '
$ws.Range("E48").Value = '  +6.02%  '
$ws.Range("E49").Value = '  -15.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.91'
$ws.Range("E50").Value = '  -14.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.578'
$ws.Range("E51").Value = '  -14.36%  '
